# Weekly update: insert the newest price-report row at the top of the
# data block (row 52), pushing all existing data rows down by one.
# Final table grows from A1:R100 to A1:R101.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 52; this shifts former rows 52..100 down
# to become rows 53..101 (EntireRow.Insert default behaviour).
$ws.Rows.Item(52).Insert()

# Populate the newly inserted row 52 with the latest week's data.
$ws.Range("A52").Value = 5
$ws.Range("B52").Value = "Macroferia Regional de Talca"
$ws.Range("C52").Value = "Maule"
$ws.Range("D52").Value = 44447
$ws.Range("E52").Value = 7
$ws.Range("F52").Value = 100112017
$ws.Range("G52").Value = "Apio"
$ws.Range("H52").Value = "Americana (o)"
$ws.Range("I52").Value = "Primera"
$ws.Range("J52").Value = 800
$ws.Range("K52").Value = 7500
$ws.Range("L52").Value = 7500
$ws.Range("M52").Value = 7500
$ws.Range("N52").Value = "`$/docena de matas"
$ws.Range("O52").Value = "Provincia del Elquí"
$ws.Range("P52").Value = 1250
$ws.Range("Q52").Value = 6
$ws.Range("R52").Value = "Hortaliza"
